$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the 8d2120e7... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-21 08:51:44"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the 8d2120e7... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-21 08:51:39"
$wsZhCn.Range("K4").Value = "2016-08-21 08:51:55"

# de-de sheet: Correspond Handoff Datetime (shares the same underlying
# "Latest HO Xliff Generate Date" text as Overview!G4) and
# Correspond Handback DateTime for the 8d2120e7... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-21 08:51:44"
$wsDeDe.Range("K4").Value = "2016-08-21 08:52:05"
